# Generate Report for Handoff
#
# The existing tracked file (0804cc9e-...md) was regenerated under a new
# temp id (7c40ef38-...md) -- its row is updated in place on every sheet.
# A second, content-duplicate file (ffff4bbc167e-...md) is handed off for
# the first time, so a brand new row is appended on every sheet.

$wb = $excel.ActiveWorkbook

$oldName = "0804cc9e-3e7a-46d5-ad3e-5ab4dd460376.md"
$newName = "7c40ef38-762e-4b7c-bd4c-740a9e180938.md"
$dupName = "ffff4bbc167e-edc6-4ca4-a78c-8a3c9e3b3e60.md"

$zhXlf = "7c40ef38-762e-4b7c-bd4c-740a9e180938.777d374e8f40e1d97e2cac394db03e87b421b0a5.zh-cn.xlf"
$deXlf = "7c40ef38-762e-4b7c-bd4c-740a9e180938.777d374e8f40e1d97e2cac394db03e87b421b0a5.de-de.xlf"

$genDate  = "2016-08-31 08:05:47"
$zhHoDate = "2016-08-31 08:05:35"
$deHoDate = "2016-08-31 08:05:47"

$repoBase = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/90b00ebe2896f20eaf6c2d1143eb0d27d3c7c56c/e2e/"

# ---------------------------------------------------------------------------
# Sheet "Overview"
# ---------------------------------------------------------------------------
$wsOv = $wb.Worksheets.Item("Overview")

# Existing row (row 2): refresh file name / path / generate date.
$wsOv.Range("B2").Hyperlinks.Delete()
$wsOv.Range("A2").Value = $newName
$wsOv.Range("B2").Value = "e2e\" + $newName
$wsOv.Range("G2").Value = $genDate
$wsOv.Hyperlinks.Add($wsOv.Range("B2"), $repoBase + $newName, "", "", "e2e\" + $newName) | Out-Null

# New row (row 3) for the duplicate file.
$wsOv.Range("A3").Value = $dupName
$wsOv.Range("B3").Value = "e2e\" + $dupName
$wsOv.Range("C3").Value = ".md"
$wsOv.Range("D2").Copy($wsOv.Range("D3"))
$wsOv.Range("E3").Value = "Ready for handoff"
$wsOv.Range("F3").Value = "Ready for handoff"
$wsOv.Range("G3").Value = $genDate
$wsOv.Hyperlinks.Add($wsOv.Range("B3"), $repoBase + $dupName, "", "", "e2e\" + $dupName) | Out-Null

# Grow the "Overview" table to include the new row.
$loOv = $wsOv.ListObjects.Item(1)
$loOv.Resize($wsOv.Range("A1:G3"))

# ---------------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

# Existing row (row 2): refresh file name + handoff xliff + datetime.
$wsZh.Range("A2").Hyperlinks.Delete()
$wsZh.Range("A2").Value = $newName
$wsZh.Range("G2").Value = $zhXlf
$wsZh.Range("H2").Value = $zhHoDate
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $repoBase + $newName, "", "", $newName) | Out-Null

# New row (row 3) for the duplicate file.
$wsZh.Range("A3").Value = $dupName
$wsZh.Range("B3").Value = ".md"
$wsZh.Range("C3").Value = "Ready for handoff"
$wsZh.Range("D3").Value = "e2e"
$wsZh.Range("E3").Value = "ht"
$wsZh.Range("F3").Value = "True"
$wsZh.Range("G3").Value = $zhXlf
$wsZh.Range("H3").Value = $zhHoDate
$wsZh.Range("I2").Copy($wsZh.Range("I3"))
$wsZh.Range("J2").Copy($wsZh.Range("J3"))
$wsZh.Range("K3").Value = "0001-01-01 00:00:00"
$wsZh.Range("L2").Copy($wsZh.Range("L3"))
$wsZh.Range("M3").Value = "True"
$wsZh.Range("N2").Copy($wsZh.Range("N3"))
$wsZh.Range("O3").Value = "False"
$wsZh.Range("P2").Copy($wsZh.Range("P3"))
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $repoBase + $dupName, "", "", $dupName) | Out-Null

$loZh = $wsZh.ListObjects.Item(1)
$loZh.Resize($wsZh.Range("A1:P3"))

# ---------------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

# Existing row (row 2): refresh file name + handoff xliff + datetime.
$wsDe.Range("A2").Hyperlinks.Delete()
$wsDe.Range("A2").Value = $newName
$wsDe.Range("G2").Value = $deXlf
$wsDe.Range("H2").Value = $deHoDate
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $repoBase + $newName, "", "", $newName) | Out-Null

# New row (row 3) for the duplicate file.
$wsDe.Range("A3").Value = $dupName
$wsDe.Range("B3").Value = ".md"
$wsDe.Range("C3").Value = "Ready for handoff"
$wsDe.Range("D3").Value = "e2e"
$wsDe.Range("E3").Value = "ht"
$wsDe.Range("F3").Value = "True"
$wsDe.Range("G3").Value = $deXlf
$wsDe.Range("H3").Value = $deHoDate
$wsDe.Range("I2").Copy($wsDe.Range("I3"))
$wsDe.Range("J2").Copy($wsDe.Range("J3"))
$wsDe.Range("K3").Value = "0001-01-01 00:00:00"
$wsDe.Range("L2").Copy($wsDe.Range("L3"))
$wsDe.Range("M3").Value = "True"
$wsDe.Range("N2").Copy($wsDe.Range("N3"))
$wsDe.Range("O3").Value = "False"
$wsDe.Range("P2").Copy($wsDe.Range("P3"))
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $repoBase + $dupName, "", "", $dupName) | Out-Null

$loDe = $wsDe.ListObjects.Item(1)
$loDe.Resize($wsDe.Range("A1:P3"))

Write-Host "Handoff report regenerated: row for $newName refreshed, row for $dupName appended on all sheets."
